# Nalco aluminium-ingot price sheet: a new day's price row is published,
# pushing the whole history table down by one row (row 2 becomes row 3,
# ..., the old last row 133 becomes row 134, i.e. gets duplicated so the
# table grows by one entry at the top).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room: insert a blank row at row 2, shifting rows 2-133 down
#        to 3-134 (row 134 automatically ends up a copy of old row 133). ---
$ws.Rows.Item(2).Insert()

# The freshly inserted row has no number formatting of its own yet; copy
# the formats (styles / column widths aside) from the row right below it
# (the old "16-12-2025" row, now row 3) so the new row renders the same
# way as every other data row.
$ws.Range("A3:F3").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122)

# --- 2. Fill in the new top row's data (latest publication date). ---
# Columns A and E hold dd-mm-yyyy strings stored as *text*, not Excel
# dates, in this workbook - force text so Excel's autodetect doesn't
# silently reinterpret an unambiguous day/month pair as a real date.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value2 = "17-12-2025"
$ws.Range("B2").Value2 = "ALUMINIUM INGOT"
$ws.Range("C2").Value2 = "IE07"
$ws.Range("D2").Value2 = 296.05
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value2 = "02-11-2025"
$ws.Range("F2").Value2 = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf"

# --- 3. Fix up the "Circular Link" hyperlinks. ---
# Row-insert shifts cell *contents* down correctly, but this host's
# hyperlink annotations stay pinned to their original row numbers, so
# after the insert they'd point at the wrong rows (and the brand-new
# row 134 would have no hyperlink at all). Rebuild them all from the
# (now-correct) displayed URL text, which is self-referential in this
# sheet - the hyperlink target always equals the cell's own text.
$ws.Hyperlinks.Delete()
For ($r = 2; $r -le 134; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $ws.Hyperlinks.Add($cell, $cell.Value2)
}

# Adding hyperlinks stamps the default blue/underlined "Hyperlink" style
# on every cell touched; restore the sheet's plain data-cell formatting
# (matching every other non-header cell) across the whole column.
$ws.Range("C3").Copy()
$ws.Range("F2:F134").PasteSpecial(-4122)
